# edit.ps1 - applies the CORE COMPETENCIES condensation and adds a new
# TECHNICAL SKILLS section, per the target diff.

$d = $word.ActiveDocument
$bullet = [char]8226

# ---------------------------------------------------------------------
# 1) Condense the three CORE COMPETENCIES detail paragraphs into a single
#    short summary line, deleting the other two paragraphs outright.
# ---------------------------------------------------------------------

$coreHeading = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq "CORE COMPETENCIES") {
        $coreHeading = $p
    }
}

$firstCompetency = $coreHeading.Next()

$firstCompetency.Range.Text = "Product Management & Strategy " + $bullet + " Technical Product Development " + $bullet + " Platform & Infrastructure"

# Re-resolve ".Next()" fresh before each delete - once a sibling paragraph
# is deleted, any previously captured reference to the paragraph after it
# becomes stale (it now points at whatever shifted into that slot).
$firstCompetency.Next().Range.Delete()
$firstCompetency.Next().Range.Delete()

# ---------------------------------------------------------------------
# 2) Append a new "TECHNICAL SKILLS" section after the last bullet of
#    "Technical Leadership & Management" (the "Trained analytical..."
#    line), before the closing "For a more detailed..." paragraph.
# ---------------------------------------------------------------------

$lastBullet = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith($bullet + " Trained analytical and engineering staff")) {
        $lastBullet = $p
    }
}

$lastBullet.Range.InsertParagraphAfter() | Out-Null
$headingPara = $lastBullet.Next()
$headingPara.Range.Text = "TECHNICAL SKILLS"

$headingPara.Range.InsertParagraphAfter() | Out-Null
$skill1 = $headingPara.Next()
$skill1.Range.Text = "PRODUCT MANAGEMENT & STRATEGY Product Conception & Ideation; Product Architecture & Design; Product Lifecycle Management; B2B SaaS Development"

$skill1.Range.InsertParagraphAfter() | Out-Null
$skill2 = $skill1.Next()
$skill2.Range.Text = "TECHNICAL PRODUCT DEVELOPMENT Full-Stack Development; Cloud Platforms; Big Data Technologies; API Development"

$skill2.Range.InsertParagraphAfter() | Out-Null
$skill3 = $skill2.Next()
$skill3.Range.Text = "PLATFORM & INFRASTRUCTURE Multi-tenant Architecture; Data Warehousing; Geospatial Platforms; Security & Compliance"

# Apply the Heading2 style to the new section heading only after all of
# the sibling body paragraphs already exist, since InsertParagraphAfter()
# otherwise propagates the heading style onto subsequently split paragraphs.
$headingPara.Style = "Heading2"

Write-Output "edit complete"
